# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (and before "总计"),
#    formatted/filled the same way as the existing "2021-Q4" sheet.
# 2. Insert a new leading data row into the "总计" sheet summarising the new
#    quarter, shifting the previous rows down and renumbering column A.

$wb = $excel.ActiveWorkbook

$q4 = $wb.Worksheets.Item("2021-Q4")

# --- 1. New "2022-Q1" sheet ------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Fetch the "总计" sheet only AFTER the insert above so it resolves to the
# actual (shifted) worksheet rather than a now-stale positional handle.
$zj = $wb.Worksheets.Item("总计")

# Copy header row + first data row formatting (font/border/alignment) from
# the "2021-Q4" sheet so the new sheet looks identical in style.
$q4.Range("B1:H1").Copy($q1.Range("B1:H1"))
$q4.Range("A2:H2").Copy($q1.Range("A2:H2"))

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

$q1.Range("A2").Value = 0
# Leading apostrophe keeps these numeric-looking values stored as TEXT
# (matching the source data, e.g. fund code "512590" / ratios as text),
# the same way the other quarter sheets store them.
$q1.Range("B2").Value = "'512590"
$q1.Range("C2").Value = "浦银安盛中证高股息精选ETF"
$q1.Range("D2").Value = "'0.59"
$q1.Range("E2").Value = "'96.43"
$q1.Range("F2").Value = "'2.18"
$q1.Range("G2").Value = "'0.0129"
$q1.Range("H2").Value = 7

# --- 2. Update "总计" sheet with the new quarter as its first data row ----
$zj.Rows("2:2").Insert()

# Row 2 (new) keeps the same look as the row that used to be there (now row 3):
# style A2 off A3 (bold/centered/bordered), clear the inherited row format on
# B2:D2 so it matches the plain look of the other data rows.
$zj.Range("A3").Copy($zj.Range("A2"))
$zj.Range("B2:D2").ClearFormats()

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 1
$zj.Range("D2").Value = 0.01

$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
